# The clue-layout grid used the lowercase letter "w" as a filler/placeholder
# value throughout the board. This pass normalizes every one of those cells
# to the uppercase "W" so it reads consistently with the other single-letter
# codes (C, D, A, M, G, B, L, K, O, X) already used on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Whole-cell (not partial/substring) replace of "w" -> "W" across every used
# cell on the sheet. None of the other codes on this sheet contain a "w" as
# a substring, so this only ever touches cells whose entire content is "w".
$used = $ws.UsedRange
[void]$used.Replace("w", "W", 1, 1, $false, $false, $false, $false)

# Leave the cursor where the author last clicked while reviewing the edit.
[void]$ws.Range("F7").Select()
